$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update employee data rows (first name, last name, username changed;
# middle name, photo, password, confirmPassword stay the same)
$ws.Range("A2").Value = "Bertokr"
$ws.Range("C2").Value = "Sapirtom"
$ws.Range("E2").Value = "Bertokk325"

$ws.Range("A3").Value = "Weportt"
$ws.Range("C3").Value = "Derakoll"
$ws.Range("E3").Value = "Weportt325"

$ws.Range("A4").Value = "Xeelopp"
$ws.Range("C4").Value = "Pomedorr"
$ws.Range("E4").Value = "Xeelopp325"

# Update the active cell selection
$ws.Range("E7").Select()
